# Re-theme the deck: replace the active "Integral" color palette with the
# "Office Theme" palette. The deck already ships both theme.xml parts; this
# rewrites the slide master's live ThemeColorScheme - i.e. the <a:clrScheme>
# entries inside ppt/theme/theme1.xml - to the Office RGB values.
#
# COM's RGB long is 0x00BBGGRR (R + G*256 + B*65536), and
# ThemeColorScheme.Item(1..12) walks the clrScheme in document order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink.

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme

# Office Theme colors, as (R,G,B) byte triples, in ThemeColorScheme order.
$officeColors = @(
    @(0x00, 0x00, 0x00),  # 1  dk1      000000
    @(0xFF, 0xFF, 0xFF),  # 2  lt1      FFFFFF
    @(0x44, 0x54, 0x6A),  # 3  dk2      44546A
    @(0xE7, 0xE6, 0xE6),  # 4  lt2      E7E6E6
    @(0x5B, 0x9B, 0xD5),  # 5  accent1  5B9BD5
    @(0xED, 0x7D, 0x31),  # 6  accent2  ED7D31
    @(0xA5, 0xA5, 0xA5),  # 7  accent3  A5A5A5
    @(0xFF, 0xC0, 0x00),  # 8  accent4  FFC000
    @(0x44, 0x72, 0xC4),  # 9  accent5  4472C4
    @(0x70, 0xAD, 0x47),  # 10 accent6  70AD47
    @(0x05, 0x63, 0xC1),  # 11 hlink    0563C1
    @(0x95, 0x4F, 0x72)   # 12 folHlink 954F72
)

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $rgb = $officeColors[$i]
    $r = $rgb[0]
    $g = $rgb[1]
    $b = $rgb[2]
    $colorScheme.Item($i + 1).RGB = $r + ($g * 256) + ($b * 65536)
}
